$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.497.19'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '1.617.59'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'210.73"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = "'22.80"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').Value = "'0.0612"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').Value = "'0.0886"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').Value = '1.848.64'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').Value = '1.607.41'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = "'0.549"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.32%  '
$ws.Range('D16').Value = "'65.10"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').Value = '27.478.64'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').Value = "'7.51"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').Value = "'4.28"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('D23').Value = "'10.16"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('E24').Value = '  +5.18%  '
$ws.Range('D25').Value = "'150.85"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('D27').Value = "'6.83"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.07%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = "'15.52"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.04%  '
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('D33').Value = '1.468.56'
$ws.Range('E33').Value = '  +1.73%  '
$ws.Range('D34').Value = "'3.07"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.30%  '
$ws.Range('E35').Value = '  -3.94%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').Value = "'0.951"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +6.14%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = "'0.558"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.26%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.0167"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = "'67.93"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.60%  '
$ws.Range('B43').Value = 'mCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D43').Value = "'2.48"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'0.988"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.13%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = "'2.21"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.23%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'5.26"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.90%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.758.66'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.71"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = "'86.55"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0105'
$ws.Range('E50').Value = '  -2.64%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.101"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.84%  '
